# Regenerate the localization-status report ("Generate Report for Handoff"):
#  - the file that was "Ready for handoff" (a963005b-...) has been replaced by a
#    newly generated one (f368b8f2-...), with a fresh handoff .xlf (f11f95e5...)
#  - a second source file (e0cb3ed1-...) now shows up with "Handoff transform failed"
#  - the ".localization-config" row shifts down to make room for the new row

$wb = $excel.ActiveWorkbook

$mdUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/abd7b4226307d46ea66c512b5a3ad7e323937eb2/e2e/f368b8f2-1195-4c2d-b37b-dbdd291f36e6.md"
$md2Url     = "https://github.com/OpenLocalizationTest/oltest/blob/abd7b4226307d46ea66c512b5a3ad7e323937eb2/e2e/e0cb3ed1-9654-47c1-b838-6b6c457e59b9.md"
$configUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/abd7b4226307d46ea66c512b5a3ad7e323937eb2/.localization-config"
$zhcnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f3bc612f3acd128dabfca35064e2a1ab10129df/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/f368b8f2-1195-4c2d-b37b-dbdd291f36e6.f11f95e5906ad50bc984d5b0103d06e5cc6bc4d4.zh-cn.xlf"
$dedeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/96bcafa9b6a9ec20596e27fcfea288155d59bfe0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/f368b8f2-1195-4c2d-b37b-dbdd291f36e6.f11f95e5906ad50bc984d5b0103d06e5cc6bc4d4.de-de.xlf"

$mdName    = "f368b8f2-1195-4c2d-b37b-dbdd291f36e6.md"
$md2Name   = "e0cb3ed1-9654-47c1-b838-6b6c457e59b9.md"
$configName = ".localization-config"
$zhcnXlfName = "f368b8f2-1195-4c2d-b37b-dbdd291f36e6.f11f95e5906ad50bc984d5b0103d06e5cc6bc4d4.zh-cn.xlf"
$dedeXlfName = "f368b8f2-1195-4c2d-b37b-dbdd291f36e6.f11f95e5906ad50bc984d5b0103d06e5cc6bc4d4.de-de.xlf"

$readyStatus = "Ready for handoff"
$failedStatus = "Handoff transform failed"
$ignoreStatus = "Not to be localized"

$zeroDate = "0001-01-01 00:00:00"
$zhcnHandoffDate = "2016-01-25 08:16:17"
$dedeHandoffDate = "2016-01-25 08:16:30"

function Set-Cells($ws, $values) {
    foreach ($addr in $values.Keys) {
        $ws.Range($addr).Value = $values[$addr]
    }
}

function Style-Link($ws, $addr) {
    # match the look of the pre-existing hyperlink cells (blue, underlined)
    $ws.Range($addr).Font.Underline = 2
    $ws.Range($addr).Font.Color = 15570276
}

function Style-DateCol($ws, $addr) {
    $ws.Range($addr).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewValues = [ordered]@{
    "A2" = $mdName;     "B2" = $readyStatus;  "C2" = $readyStatus
    "A3" = $md2Name;    "B3" = $failedStatus; "C3" = $failedStatus
    "A4" = $configName; "B4" = $ignoreStatus; "C4" = $ignoreStatus
}
Set-Cells $wsOverview $overviewValues

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $mdName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $md2Url, "", "", $md2Name)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $configUrl, "", "", $configName)
Style-Link $wsOverview "A2"
Style-Link $wsOverview "A3"
Style-Link $wsOverview "A4"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhValues = [ordered]@{
    "A2" = $mdName;     "B2" = $readyStatus;  "C2" = $zhcnXlfName; "D2" = $zhcnHandoffDate; "G2" = $zeroDate; "H2" = "Include"
    "A3" = $md2Name;    "B3" = $failedStatus;                      "D3" = $zeroDate;        "G3" = $zeroDate; "H3" = "Ignored"
    "A4" = $configName; "B4" = $ignoreStatus;                      "D4" = $zeroDate;        "G4" = $zeroDate; "H4" = "Ignored"
}
Set-Cells $wsZh $zhValues
Style-DateCol $wsZh "D2"
Style-DateCol $wsZh "D3"
Style-DateCol $wsZh "D4"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhcnXlfUrl, "", "", $zhcnXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $md2Url, "", "", $md2Name)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $configUrl, "", "", $configName)
Style-Link $wsZh "A2"
Style-Link $wsZh "C2"
Style-Link $wsZh "A3"
Style-Link $wsZh "A4"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deValues = [ordered]@{
    "A2" = $mdName;     "B2" = $readyStatus;  "C2" = $dedeXlfName; "D2" = $dedeHandoffDate; "G2" = $zeroDate; "H2" = "Include"
    "A3" = $md2Name;    "B3" = $failedStatus;                      "D3" = $zeroDate;        "G3" = $zeroDate; "H3" = "Ignored"
    "A4" = $configName; "B4" = $ignoreStatus;                      "D4" = $zeroDate;        "G4" = $zeroDate; "H4" = "Ignored"
}
Set-Cells $wsDe $deValues
Style-DateCol $wsDe "D2"
Style-DateCol $wsDe "D3"
Style-DateCol $wsDe "D4"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $dedeXlfUrl, "", "", $dedeXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $md2Url, "", "", $md2Name)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $configUrl, "", "", $configName)
Style-Link $wsDe "A2"
Style-Link $wsDe "C2"
Style-Link $wsDe "A3"
Style-Link $wsDe "A4"
